$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New service-log rows appended after the existing last row (392).
# Columns: A=DATE, B=VEHICLE REG NO, C=VEHICLE BRAND, D=ISSUE, E=STATUS, F=AMOUNT, G=CASH TYPE

$ws.Cells.Item(393, 1).Value = 44817
$ws.Cells.Item(393, 2).Value = "KA05MP4915"
$ws.Cells.Item(393, 3).Value = "PUNTO"
$ws.Cells.Item(393, 4).Value = "HAFT SAFT CHANGE"
$ws.Cells.Item(393, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(394, 1).Value = 44817
$ws.Cells.Item(394, 2).Value = "KA03MV0746"
$ws.Cells.Item(394, 3).Value = "ERTIGA"
$ws.Cells.Item(394, 4).Value = "PMS"
$ws.Cells.Item(394, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(395, 1).Value = 44817
$ws.Cells.Item(395, 2).Value = "KA53P6612"
$ws.Cells.Item(395, 3).Value = "RITZ"
$ws.Cells.Item(395, 4).Value = "ECM CHANGE"
$ws.Cells.Item(395, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(395, 6).Value = 9500
$ws.Cells.Item(395, 7).Value = "CARD"

$ws.Cells.Item(396, 1).Value = 44817
$ws.Cells.Item(396, 2).Value = "KA03NE7365"
$ws.Cells.Item(396, 3).Value = "NEXON"
$ws.Cells.Item(396, 4).Value = "BODY SHOP"
$ws.Cells.Item(396, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(397, 1).Value = 44817
$ws.Cells.Item(397, 2).Value = "KA03MP9135"
$ws.Cells.Item(397, 3).Value = "FIGO"
$ws.Cells.Item(397, 4).Value = "BODY SHOP"
$ws.Cells.Item(397, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(398, 1).Value = 44817
$ws.Cells.Item(398, 2).Value = "AP28D1658"
$ws.Cells.Item(398, 3).Value = "WAGON R"
$ws.Cells.Item(398, 4).Value = "BODY SHOP"
$ws.Cells.Item(398, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(399, 1).Value = 44817
$ws.Cells.Item(399, 2).Value = "KA35M1408"
$ws.Cells.Item(399, 3).Value = "ACCENT"
$ws.Cells.Item(399, 4).Value = "PMS"
$ws.Cells.Item(399, 5).Value = "WORK IN PROGRESS"

$ws.Cells.Item(400, 1).Value = 44817
$ws.Cells.Item(400, 2).Value = "KA05M6661"
$ws.Cells.Item(400, 3).Value = "I20"
$ws.Cells.Item(400, 4).Value = "BODY SHOP"
$ws.Cells.Item(400, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(400, 6).Value = 24056

# Match the author's final selection state.
$null = $ws.Range("A401").Select()
